# re-added in BLEU and GLEU score calculations
#
# NOTE: cells are written in an order chosen so that the workbook's shared
# string table is (re)built in the same order as the original commit, with
# the "Decoded sequence"/"Input sequence" repeat strings deduplicated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = ' different than what i expected on the outset but i still ended up having a good time with it eos'
$ws.Range("C3").Value  = 'Decoded sequence: sos i the the the the the the'
$ws.Range("A13").Value = 'Input sequence: sos spring 2015 plastic memories danmachi kekkai sensen only one of those i liked eos'
$ws.Range("C4").Value  = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A15").Value = 'Input sequence: sos to all the saber lovers out there enjoy https imgurcomakdrqu some content is nsfw eos'
$ws.Range("A17").Value = 'Input sequence: names are and where to start there s lots of good anime out there today some long some short eos'
$ws.Range("A18").Value = 'Decoded sequence: sos the the the the the the the this'
$ws.Range("A19").Value = 'Input sequence: sos first and last episodes of plastic memories eos'
$ws.Range("A21").Value = 'Input sequence: sos winter 2015 and i watched pretty much everything but my favorites were death parade parasyte and yuri kuma eos'
$ws.Range("A23").Value = 'Input sequence: 20140302110443 and basketball lesbians http i1294photobucketcomalbumsb619isonnazzoanimelarge zpssf3c0razgif also the only one who can kill me is me https smediacacheak0pinimgcom236xdac909dac90980a052e3bde464b4cad968e011jpg eos'
$ws.Range("A25").Value = 'Input sequence: sos i think you d enjoy psychopass http myanimelistnetanime13601psychopass eos'
$ws.Range("A27").Value = 'Input sequence: for some reason bears all over the world rose up and attacked humanity yuri kuma was a fun show eos'
$ws.Range("A29").Value = 'Input sequence: be immersed in the story sadly it s only two episodes and the novel s translation are barely breathing eos'

# Remaining cells reuse the "Decoded sequence: sos the the ..." string already
# introduced above (C4), so the shared string table stays de-duplicated.
$ws.Range("A16").Value = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A20").Value = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A22").Value = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A24").Value = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A26").Value = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A28").Value = 'Decoded sequence: sos the the the the the the the the the the the the'
$ws.Range("A30").Value = 'Decoded sequence: sos the the the the the the the the the the the the'

# Reproduce the sheet's recorded selection/active cell
$ws.Range("A15").Select() | Out-Null
